$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: morning/evening weight readings were swapped, morning body fat
# corrected, and the medication dosage entry removed.
$ws.Range("B2").Value = 104.7
$ws.Range("C2").Value = 105.1
$ws.Range("D2").Value = 30.3
$ws.Range("F2").ClearContents()

# Rows 3-7: the medication dosage entries were removed.
$ws.Range("F3").ClearContents()
$ws.Range("F4").ClearContents()
$ws.Range("F5").ClearContents()
$ws.Range("F6").ClearContents()
$ws.Range("F7").ClearContents()

# Row 8: evening weight, evening body fat, and dosage entries were removed.
$ws.Range("C8").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("F8").ClearContents()
